$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The nowcast table grows by one more run: dates shift forward one slot and a
# new row (2025-08-30) is appended; every revision figure is refreshed too.

# Column A holds date labels as text (e.g. "2025-03-30"); format as Text first
# so Excel does not auto-convert the strings into date serial values, then
# restore the default "Normal" cell style so formatting matches the rest of
# the sheet once the text has been committed.
$ws.Range("A2:A12").NumberFormat = "@"
$ws.Range("A2").Value = "2025-03-30"
$ws.Range("A3").Value = "2025-04-15"
$ws.Range("A4").Value = "2025-04-30"
$ws.Range("A5").Value = "2025-05-15"
$ws.Range("A6").Value = "2025-05-30"
$ws.Range("A7").Value = "2025-06-15"
$ws.Range("A8").Value = "2025-06-30"
$ws.Range("A9").Value = "2025-07-15"
$ws.Range("A10").Value = "2025-07-30"
$ws.Range("A11").Value = "2025-08-15"
$ws.Range("A12").Value = "2025-08-30"
$ws.Range("A2:A12").Style = "Normal"

# Updated nowcast / revision figures (columns B-K), rows 2-12
$ws.Range("B2").Value = 0.28020599571202126
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("B3").Value = 0.3015352712016185
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = -0.005305442770179054
$ws.Range("E3").Value = 0.00011983417635948429
$ws.Range("F3").Value = -0.0012382981285977287
$ws.Range("G3").Value = 0.0002294026952707526
$ws.Range("H3").Value = -0.00009233597345637557
$ws.Range("I3").Value = -0.0005033625729192658
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = -0.00006149492768747766
$ws.Range("B4").Value = 0.2949042652081544
$ws.Range("C4").Value = -0.002299089259277439
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.0005283873919141694
$ws.Range("F4").Value = 0.00007485671976047254
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.0001528094164175833
$ws.Range("I4").Value = -0.002830200919053255
$ws.Range("J4").Value = 0.0005152425631643249
$ws.Range("K4").Value = 0.00028448530148811324
$ws.Range("B5").Value = 0.2981567949247489
$ws.Range("C5").Value = 0.010802406384145464
$ws.Range("D5").Value = -0.007329500168472452
$ws.Range("E5").Value = 0.00028864675626416483
$ws.Range("F5").Value = 0.0011639600146210885
$ws.Range("G5").Value = -0.0015942340198711951
$ws.Range("H5").Value = 0.000023392597581988896
$ws.Range("I5").Value = -0.0006855881944997005
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0.000057677842711600924
$ws.Range("B6").Value = 0.23812082277261054
$ws.Range("C6").Value = 0.024633384853765677
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = -0.0003166319476319479
$ws.Range("F6").Value = 0.00002152181581803508
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = -0.00009403146723937032
$ws.Range("I6").Value = -0.0024086954913735653
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0.0000101858813182365
$ws.Range("B7").Value = 0.24666497396253714
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = -0.0029127972896584313
$ws.Range("E7").Value = -0.0014236081213807384
$ws.Range("F7").Value = -0.007530013574321624
$ws.Range("G7").Value = 0.0012777311053934048
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0.0003195201616695691
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0.0006694245554834388
$ws.Range("B8").Value = 0.2139665852334336
$ws.Range("C8").Value = -0.05903436695833248
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = -0.000043426523151359026
$ws.Range("F8").Value = -0.00022989267257996694
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0.00006080320142694691
$ws.Range("I8").Value = 0.0009194365758756895
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0.0011855328643549934
$ws.Range("B9").Value = 0.2350475647919402
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0.006308588744644696
$ws.Range("E9").Value = -0.0031762827858649286
$ws.Range("F9").Value = -0.004269166799258723
$ws.Range("G9").Value = 0.002259503082378708
$ws.Range("H9").Value = -0.000193213504180704
$ws.Range("I9").Value = 0.00011021588056233323
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0.00005887571740167963
$ws.Range("B10").Value = 0.3027777379141217
$ws.Range("C10").Value = 0.09986515018166023
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = -0.0008579118981499199
$ws.Range("F10").Value = -0.00015095194325885408
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = -0.00001368701595864286
$ws.Range("I10").Value = 0.0008278385484764119
$ws.Range("J10").Value = -0.002209778723577333
$ws.Range("K10").Value = -0.000026371346599218803
$ws.Range("B11").Value = 0.35843433538826414
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = -0.03372662557929491
$ws.Range("E11").Value = 0.0031890477873551537
$ws.Range("F11").Value = 0.005892803768770394
$ws.Range("G11").Value = 0.003463564270857285
$ws.Range("H11").Value = 0.0010047848573160613
$ws.Range("I11").Value = 0.0037498874306605636
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0.0013848745584090216
$ws.Range("B12").Value = 0.24032812944156748
$ws.Range("C12").Value = -0.07004065415646896
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0.0024257726423090385
$ws.Range("F12").Value = 0.00004083159807914327
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0.000002780460412795232
$ws.Range("I12").Value = -0.002801301943813378
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0.00017432110062848283
